$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add headers for the new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header style used by the other header cells (bold, centered, thin border)
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Fill in the team record (Wins/Losses/Ties) for every data row (2 through 44)
for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 30).Value = 102
    $ws.Cells.Item($row, 31).Value = 60
    $ws.Cells.Item($row, 32).Value = 0
}
